# Improve data parsing logic
# Appends a new row 26 (mirroring the layout of the existing data rows)
# to each of the four worksheets, using the latest parsed packet values.

$wb = $excel.ActiveWorkbook

$rowData = @{
    1 = @{
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x8C"
        E = "0x07"
        F = 400.0
        G = 568631262647113000000000.0
        H = 396.0
        I = 7.0
    }
    2 = @{
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x78"
        E = "0x19"
        F = 380.0
        G = 568432987514711000000000.0
        H = 376.0
        I = 25.0
    }
    3 = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x6D"
        E = "0x15"
        F = 110.0
        G = 568631262647113000000000.0
        H = 109.0
        I = 15.0
    }
    4 = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x82"
        E = "0x9"
        F = 130.0
        G = 568631262647113000000000.0
        H = 130.0
        I = 9.0
    }
}

$newTimestamp = 45812.46581018518

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowData[$i]
    $newRow = 26

    $ws.Range("A" + $newRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A" + $newRow).Value = $newTimestamp

    $ws.Range("B" + $newRow).Value = $data.B
    $ws.Range("C" + $newRow).Value = $data.C
    $ws.Range("D" + $newRow).Value = $data.D
    $ws.Range("E" + $newRow).Value = $data.E
    $ws.Range("F" + $newRow).Value = $data.F
    $ws.Range("G" + $newRow).Value = $data.G
    $ws.Range("H" + $newRow).Value = $data.H
    $ws.Range("I" + $newRow).Value = $data.I
}
